# TournRPG-115: 全体攻撃の実装 (implement all-enemy/all-ally attack range)
# Rename the Japanese "range" values on the skill sheet to their English
# equivalents, and widen column F ("range") so the longer English labels
# are not clipped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("skill")

# F3/F4: 敵単 -> ENEMY_ONE   (shared by SKILL001 二段攻撃 rows 3 & 4)
$ws.Range("F3").Value = "ENEMY_ONE"
$ws.Range("F4").Value = "ENEMY_ONE"

# F5: 敵全 -> ENEMY_ALL      (SKILL003 毒の粉)
$ws.Range("F5").Value = "ENEMY_ALL"

# F6: 味単 -> FRIEND_ONE     (SKILL004 HP回復)
$ws.Range("F6").Value = "FRIEND_ONE"

# F7: 味全 -> FRIEND_ALL     (SKILL005 攻撃アップ)
$ws.Range("F7").Value = "FRIEND_ALL"

# Widen column F (range) from 4.75 to ~8.75 chars to fit the new English labels
$ws.Columns("F").ColumnWidth = 8
